$d = $word.ActiveDocument
$p = $d.Paragraphs(1)

# Add a 5-point-space paragraph border (top/left/bottom/right) to the first paragraph.
$borders = $p.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p.Range.ParagraphFormat.LeftIndent = 11.25

# Update the placeholder text and collapse the trailing-space run into it.
$d.Content.Find.Execute("**ID__AFFARS_5328_topic_7__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_SUBPART_5328_3__ID**", 2)
